# Add the new row of data for date 12/9/2025 (serial 46000) to the
# "Daily 100 Error Counts" sheet, and update the active selection to
# reflect the newly entered row, as in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = 46000
$ws.Range("B43").Value = 649
$ws.Range("C43").Value = 16
$ws.Range("D43").Value = 633

$ws.Range("A43:D43").Select() | Out-Null
